$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 69
$ws.Range('A69').Value = 111785190
$ws.Range('B69').Value = 94134
$ws.Range('D69').Value = 'NT'
$ws.Range('E69').Value = 53
$ws.Range('F69').Value = 'Vedtrappmossa'
$ws.Range('G69').Value = 'Crossocalyx hellerianus'
$ws.Range('H69').Value = '(Nees ex Lindenb.) Meyl.'
$ws.Range('Q69').Value = 577242.972410051
$ws.Range('R69').Value = 6944540.945152098

# Row 70
$ws.Range('A70').Value = 111785230
$ws.Range('B70').Value = 78578
$ws.Range('D70').Value = 'NT'
$ws.Range('E70').Value = 6458
$ws.Range('F70').Value = 'Lunglav'
$ws.Range('G70').Value = 'Lobaria pulmonaria'
$ws.Range('H70').Value = '(L.) Hoffm.'
$ws.Range('Q70').Value = 577261.8704127767
$ws.Range('R70').Value = 6944620.109213427

# Row 71
$ws.Range('A71').Value = 111785244
$ws.Range('B71').Value = 96348
$ws.Range('D71').Value = 'VU'
$ws.Range('E71').Value = 220787
$ws.Range('F71').Value = 'Knärot'
$ws.Range('G71').Value = 'Goodyera repens'
$ws.Range('H71').Value = '(L.) R. Br.'
$ws.Range('Q71').Value = 577364.1822193989
$ws.Range('R71').Value = 6944621.578847388

# Row 72
$ws.Range('A72').Value = 111785206
$ws.Range('B72').Value = 77268
$ws.Range('D72').Value = 'NT'
$ws.Range('E72').Value = 228912
$ws.Range('F72').Value = 'Mörk kolflarnlav'
$ws.Range('G72').Value = 'Carbonicola myrmecina'
$ws.Range('H72').Value = '(Ach.) Bendiksby & Timdal'
$ws.Range('Q72').Value = 577235.6798241453
$ws.Range('R72').Value = 6944655.86623876

# Row 73
$ws.Range('A73').Value = 111785202
$ws.Range('B73').Value = 78512
$ws.Range('D73').Value = 'LC'
$ws.Range('E73').Value = 6456
$ws.Range('F73').Value = 'Skinnlav'
$ws.Range('G73').Value = 'Leptogium saturninum'
$ws.Range('H73').Value = '(Dicks.) Nyl.'
$ws.Range('Q73').Value = 577215.0430418774
$ws.Range('R73').Value = 6944631.445974576

# Row 74
$ws.Range('A74').Value = 111785200
$ws.Range('B74').Value = 78512
$ws.Range('D74').Value = 'LC'
$ws.Range('E74').Value = 6456
$ws.Range('F74').Value = 'Skinnlav'
$ws.Range('G74').Value = 'Leptogium saturninum'
$ws.Range('H74').Value = '(Dicks.) Nyl.'
$ws.Range('Q74').Value = 577256.110519147
$ws.Range('R74').Value = 6944531.123615563

# Row 75
$ws.Range('A75').Value = 111785235
$ws.Range('B75').Value = 77267
$ws.Range('D75').Value = 'NT'
$ws.Range('E75').Value = 6446
$ws.Range('F75').Value = 'Kolflarnlav'
$ws.Range('G75').Value = 'Carbonicola anthracophila'
$ws.Range('H75').Value = '(Nyl.) Bendiksby & Timdal'
$ws.Range('Q75').Value = 577226.625646919
$ws.Range('R75').Value = 6944648.749557905

# Row 76
$ws.Range('A76').Value = 111785201
$ws.Range('B76').Value = 78512
$ws.Range('D76').Value = 'LC'
$ws.Range('E76').Value = 6456
$ws.Range('F76').Value = 'Skinnlav'
$ws.Range('G76').Value = 'Leptogium saturninum'
$ws.Range('H76').Value = '(Dicks.) Nyl.'
$ws.Range('Q76').Value = 577248.2772659193
$ws.Range('R76').Value = 6944530.940753835

# Row 77
$ws.Range('A77').Value = 111785229
$ws.Range('B77').Value = 78578
$ws.Range('D77').Value = 'NT'
$ws.Range('E77').Value = 6458
$ws.Range('F77').Value = 'Lunglav'
$ws.Range('G77').Value = 'Lobaria pulmonaria'
$ws.Range('H77').Value = '(L.) Hoffm.'
$ws.Range('Q77').Value = 577208.3826684169
$ws.Range('R77').Value = 6944521.722980071

# Row 78
$ws.Range('A78').Value = 111785191
$ws.Range('B78').Value = 89405
$ws.Range('D78').Value = 'NT'
$ws.Range('E78').Value = 1202
$ws.Range('F78').Value = 'Ullticka'
$ws.Range('G78').Value = 'Phellinidium ferrugineofuscum'
$ws.Range('H78').Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range('Q78').Value = 577235.6798241453
$ws.Range('R78').Value = 6944655.86623876

# Row 79
$ws.Range('A79').Value = 111785251
$ws.Range('B79').Value = 93161
$ws.Range('D79').Value = 'VU'
$ws.Range('E79').Value = 1079
$ws.Range('F79').Value = 'Aspfjädermossa'
$ws.Range('G79').Value = 'Neckera pennata'
$ws.Range('H79').Value = 'Hedw.'
$ws.Range('Q79').Value = 577283.2535308318
$ws.Range('R79').Value = 6944533.598891968

# Row 80
$ws.Range('A80').Value = 111785228
$ws.Range('B80').Value = 78578
$ws.Range('D80').Value = 'NT'
$ws.Range('E80').Value = 6458
$ws.Range('F80').Value = 'Lunglav'
$ws.Range('G80').Value = 'Lobaria pulmonaria'
$ws.Range('H80').Value = '(L.) Hoffm.'
$ws.Range('Q80').Value = 577256.110519147
$ws.Range('R80').Value = 6944531.123615563

# Row 81
$ws.Range('A81').Value = 111785192
$ws.Range('B81').Value = 89405
$ws.Range('D81').Value = 'NT'
$ws.Range('E81').Value = 1202
$ws.Range('F81').Value = 'Ullticka'
$ws.Range('G81').Value = 'Phellinidium ferrugineofuscum'
$ws.Range('H81').Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range('Q81').Value = 577281.7951240344
$ws.Range('R81').Value = 6944714.487089146

# Row 82
$ws.Range('A82').Value = 111785199
$ws.Range('B82').Value = 89416
$ws.Range('D82').Value = 'LC'
$ws.Range('E82').Value = 1205
$ws.Range('F82').Value = 'Stor aspticka'
$ws.Range('G82').Value = 'Phellinus populicola'
$ws.Range('H82').Value = 'Niemelä'
$ws.Range('Q82').Value = 577256.110519147
$ws.Range('R82').Value = 6944531.123615563
